$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# --- Row 13: the previously-blank cells get filled in with the literal
#     text "nan" (L13, N13, O13 already hold values and stay unchanged) ---
$row13Cols = "B","C","D","E","F","G","H","I","J","K","M","P","Q","R"
foreach ($col in $row13Cols) {
    $ws.Range($col + "13").Value = "nan"
}

# --- Row 14: brand-new service-card event row ---
# Card number, stored as text (matches the rest of column A).
$ws.Range("A14").Value = "'15"
$ws.Range("A14").Style = "Normal"

$ws.Range("L14").Value = "12\12\2024"
$ws.Range("N14").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O14").Value = "تيم العمل"

# Remaining row-14 cells exist but are blank text cells (force them into
# existence as empty text via the quote-prefix trick, then drop the
# quote-prefix style it leaves behind).
$row14BlankCols = "B","C","D","E","F","G","H","I","J","K","M","P","Q","R"
foreach ($col in $row14BlankCols) {
    $addr = $col + "14"
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}
